# Updates cryptos list figures (prices / 1h volume %) and swaps the
# Toncoin / InjectiveProtocol rows (28 and 29), matching the GitHub Actions
# data refresh commit.
#
# Every assigned value is prefixed with a leading apostrophe so Excel stores
# it as literal text (the sheet's Price/Volume columns hold plain text, not
# numbers - e.g. "305.77" must stay the string "305.77", not be coerced into
# a numeric cell). Resetting Style to 'Normal' immediately afterwards clears
# the quote-prefix cell style that the apostrophe trick leaves behind, so the
# cell's style index is left exactly as it was before (no style changes).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''44.000.32'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '''  -0.08%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = '''2.237.02'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '''  -0.71%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = '''  +0.17%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = '''305.77'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '''  -4.35%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = '''95.27'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '''  -6.31%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = '''0.569'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '''  -1.38%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('E8').Value = '''  +0.20%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = '''0.524'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '''  -5.18%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = '''34.57'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '''  -7.56%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = '''0.0808'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '''  -2.83%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = '''7.20'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '''  -5.06%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = '''0.104'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '''  -2.66%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = '''2.576.60'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '''  -0.81%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = '''2.238.26'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '''  -0.66%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = '''0.819'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '''  -4.26%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = '''13.50'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '''  -6.39%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = '''43.848.08'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '''  -0.15%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = '''0.0₃0957'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '''  -2.80%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = '''12.27'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '''  -7.35%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = '''6.20'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '''  -4.71%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = '''64.75'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '''  -1.47%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = '''237.46'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '''  +0.92%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = '''2.91'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '''  -7.58%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('E25').Value = '''  +0.39%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('E26').Value = '''  -7.65%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = '''9.86'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '''  -3.57%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('B28').Value = '''InjectiveProtocol'
$ws.Range('B28').Style = 'Normal'
$ws.Range('C28').Value = '''https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('C28').Style = 'Normal'
$ws.Range('D28').Value = '''36.87'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '''  -2.71%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('B29').Value = '''Toncoin'
$ws.Range('B29').Style = 'Normal'
$ws.Range('C29').Value = '''https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('C29').Style = 'Normal'
$ws.Range('D29').Value = '''2.13'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '''  -0.86%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = '''20.07'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '''  -0.57%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = '''5.89'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '''  -4.82%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = '''154.62'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '''  -4.48%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = '''0.0807'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '''  -5.16%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = '''3.31'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '''  +9.91%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('E35').Value = '''  -2.50%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = '''0.110'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '''  -4.63%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('E37').Value = '''  -0.59%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('E38').Value = '''  -7.90%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = '''15.13'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '''  -10.59%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = '''3.36'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '''  -9.70%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = '''3.79'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '''  -9.96%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('E42').Value = '''  -4.92%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('D44').Value = '''1.743.75'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '''  -2.68%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = '''84.45'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '''  +2.70%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('E46').Value = '''  -6.25%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = '''99.81'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '''  -4.70%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = '''4.90'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '''  -5.65%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = '''14.68'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '''  -1.81%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = '''8.07'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '''  -3.51%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = '''68.65'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '''  -9.08%  '
$ws.Range('E51').Style = 'Normal'
